$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the full target range to Text format so numeric- and date-looking
# strings (e.g. "5/29/2008", "173900.0") are preserved verbatim as text,
# matching the source CSV import behavior (all cells are inline/text strings).
$ws.Range("A1:P6").NumberFormat = "@"

$ws.Range("A1").Value = 'Servicer ID'
$ws.Range("B1").Value = 'Lender Loan Number'
$ws.Range("C1").Value = 'Certificate Number'
$ws.Range("D1").Value = 'Borrower Name'
$ws.Range("E1").Value = 'Property Address'
$ws.Range("F1").Value = 'Original Effective Date'
$ws.Range("G1").Value = 'Loan Amount'
$ws.Range("H1").Value = 'Billing Date'
$ws.Range("I1").Value = 'Invoice Number'
$ws.Range("J1").Value = 'Current Principal Balance'
$ws.Range("K1").Value = 'Billing Rate'
$ws.Range("L1").Value = 'Premium Due'
$ws.Range("M1").Value = 'Tax Rate'
$ws.Range("N1").Value = 'Tax Amount'
$ws.Range("O1").Value = 'Total Amount'
$ws.Range("P1").Value = 'MI Paid By'
$ws.Range("A2").Value = '16812-0001-0'
$ws.Range("B2").Value = '315232-437'
$ws.Range("C2").Value = '19732096'
$ws.Range("D2").Value = 'Aherne, Brian'
$ws.Range("E2").Value = '199 Pine Grove Drive - Pittsfield, MA 01201-0000'
$ws.Range("F2").Value = '5/29/2008'
$ws.Range("G2").Value = '173900.0'
$ws.Range("H2").Value = '5/29/2020'
$ws.Range("I2").Value = '340064'
$ws.Range("J2").Value = '173900.0'
$ws.Range("K2").Value = '0.17'
$ws.Range("L2").Value = '24.64'
$ws.Range("M2").Value = '0.0'
$ws.Range("N2").Value = '0.0'
$ws.Range("O2").Value = '24.64'
$ws.Range("P2").Value = 'Borrower'
$ws.Range("A3").Value = '16812-0001-0'
$ws.Range("B3").Value = '315232-437'
$ws.Range("C3").Value = '19732096'
$ws.Range("D3").Value = 'Aherne, Brian'
$ws.Range("E3").Value = '199 Pine Grove Drive - Pittsfield, MA 01201-0000'
$ws.Range("F3").Value = '5/29/2008'
$ws.Range("G3").Value = '173900.0'
$ws.Range("H3").Value = '6/29/2020'
$ws.Range("I3").Value = '342274'
$ws.Range("J3").Value = '173900.0'
$ws.Range("K3").Value = '0.17'
$ws.Range("L3").Value = '24.64'
$ws.Range("M3").Value = '0.0'
$ws.Range("N3").Value = '0.0'
$ws.Range("O3").Value = '24.64'
$ws.Range("P3").Value = 'Borrower'
$ws.Range("A4").Value = '16812-0001-0'
$ws.Range("B4").Value = '315232-437'
$ws.Range("C4").Value = '19732096'
$ws.Range("D4").Value = 'Aherne, Brian'
$ws.Range("E4").Value = '199 Pine Grove Drive - Pittsfield, MA 01201-0000'
$ws.Range("F4").Value = '5/29/2008'
$ws.Range("G4").Value = '173900.0'
$ws.Range("H4").Value = '7/29/2020'
$ws.Range("I4").Value = '344454'
$ws.Range("J4").Value = '173900.0'
$ws.Range("K4").Value = '0.17'
$ws.Range("L4").Value = '24.64'
$ws.Range("M4").Value = '0.0'
$ws.Range("N4").Value = '0.0'
$ws.Range("O4").Value = '24.64'
$ws.Range("P4").Value = 'Borrower'
$ws.Range("A5").Value = '16812-0001-0'
$ws.Range("B5").Value = '7602359142'
$ws.Range("C5").Value = '68441331'
$ws.Range("D5").Value = 'Gross, Tracy'
$ws.Range("E5").Value = '1208 Williamsburg Drive - Champaign, IL 61821-0000'
$ws.Range("F5").Value = '9/1/2012'
$ws.Range("G5").Value = '79705.0'
$ws.Range("H5").Value = '6/1/2020'
$ws.Range("I5").Value = '342274'
$ws.Range("J5").Value = '71371.0'
$ws.Range("K5").Value = '0.64'
$ws.Range("L5").Value = '42.51'
$ws.Range("M5").Value = '0.0'
$ws.Range("N5").Value = '0.0'
$ws.Range("O5").Value = '42.51'
$ws.Range("P5").Value = 'Borrower'
$ws.Range("A6").Value = '16812-0001-0'
$ws.Range("B6").Value = '7602359142'
$ws.Range("C6").Value = '68441331'
$ws.Range("D6").Value = 'Gross, Tracy'
$ws.Range("E6").Value = '1208 Williamsburg Drive - Champaign, IL 61821-0000'
$ws.Range("F6").Value = '9/1/2012'
$ws.Range("G6").Value = '79705.0'
$ws.Range("H6").Value = '7/1/2020'
$ws.Range("I6").Value = '344454'
$ws.Range("J6").Value = '71371.0'
$ws.Range("K6").Value = '0.64'
$ws.Range("L6").Value = '42.51'
$ws.Range("M6").Value = '0.0'
$ws.Range("N6").Value = '0.0'
$ws.Range("O6").Value = '42.51'
$ws.Range("P6").Value = 'Borrower'
